$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Date "
$ws.Range("B1").Value = "Time Spent"
$ws.Range("C1").Value = "Topics Covered"
$ws.Range("D1").Value = "File name "
$ws.Range("E1").Value = "Remarks"
$ws.Range("F1").Value = "Certifications "

# Data row
$ws.Range("A2").Value = 45000
$ws.Range("A2").NumberFormat = "mm-dd-yy"

$ws.Range("B2").Value = "3-4 Hours"
$ws.Range("C2").Value = "Git & Github- Branching, Merging, Rebase and other basic commands"
$ws.Range("D2").Value = "1)Installation and basic commands  2) Branching  3) auto email  4) Tags         5) Merging vs Rebase"
$ws.Range("E2").Value = "Completed"
$ws.Range("F2").Value = "https://www.udemy.com/certificate/UC-a8fc2a6b-4f1b-4ce1-9af5-921adb6fd73d/"

# Wrap text for the long-content cells
$ws.Range("C2").WrapText = $true
$ws.Range("D2").WrapText = $true
$ws.Range("F2").WrapText = $true

# Column widths (approximate best-fit widths from the original file)
$ws.Columns.Item(1).ColumnWidth = 9.416666666666668
$ws.Columns.Item(2).ColumnWidth = 9.083333333333334
$ws.Columns.Item(3).ColumnWidth = 12.583333333333334
$ws.Columns.Item(4).ColumnWidth = 10.583333333333334
$ws.Columns.Item(6).ColumnWidth = 11.416666666666668

# Row height for the wrapped second row
$ws.Rows.Item(2).RowHeight = 129.6

# Restore the active selection
$ws.Range("G5").Select() | Out-Null
